$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Rewards -> Steps On Failure
$ws.Range("A2").Value = "Steps On Failure"
$ws.Range("B2").Value = "no failures!"
$ws.Range("C2").Value = "no failures!"
$ws.Range("D2").Value = 569.3
$ws.Range("E2").Value = 379.01

# Row 3: Steps -> Steps On Success
$ws.Range("A3").Value = "Steps On Success"
$ws.Range("B3").Value = 127.7549019607843
$ws.Range("C3").Value = 126.025
$ws.Range("D3").Value = 277.6
$ws.Range("E3").Value = 276.15
